# Update metric values in the worksheet to reflect refreshed model run results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (MSE) updates
$ws.Cells.Item(2, 2).Value = 0.4108830721519875
$ws.Cells.Item(3, 2).Value = 0.1647362719327807
$ws.Cells.Item(4, 2).Value = 0.2192987281846223
$ws.Cells.Item(5, 2).Value = 0.3638702225807679
$ws.Cells.Item(7, 2).Value = 0.07514644587374564
$ws.Cells.Item(8, 2).Value = 0.04215534119371416
$ws.Cells.Item(9, 2).Value = 0.07796894984218636

# Column D (MAE) updates
$ws.Cells.Item(2, 4).Value = 0.525669270734418
$ws.Cells.Item(4, 4).Value = 0.3777406617731507
$ws.Cells.Item(7, 4).Value = 0.211919863475561
$ws.Cells.Item(8, 4).Value = 0.1361288253571671
$ws.Cells.Item(9, 4).Value = 0.1911874935925046

# Column G (Elapsed Time) and Column H (CPU) updates for all data rows (2-14)
for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 7).Value = 1.558460351833249
    $ws.Cells.Item($row, 8).Value = 0.9990000000000001
}
